$d = $word.ActiveDocument

function New-OpenXmlPackage([string]$bodyFragment) {
    $pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
    $pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $pkgHeader + '<w:body>' + $bodyFragment + '</w:body>' + $pkgFooter
}

function Set-ParagraphXml([int]$paraIndex, [string]$newParaXml) {
    $para = $d.Paragraphs($paraIndex)
    $rng = $para.Range
    # Exclude the trailing paragraph-mark character so the paragraph's own
    # end-of-paragraph mark (and thus its pPr) is preserved by the host
    # paragraph; we supply the replacement pPr explicitly inside $newParaXml.
    $target = $d.Range($rng.Start, $rng.End - 1)
    $xmlPkg = New-OpenXmlPackage $newParaXml
    $target.InsertXML($xmlPkg)
}

# --- Paragraph: "El sistema consulta la base de datos ..." (Flujo normal, step 1)
$p1 = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
      '<w:r><w:t xml:space="preserve">El sistema consulta la base de datos y recupera la información de la </w:t></w:r>' +
      '<w:r><w:t>VENTA (</w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>NoVenta</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>FechaRegistro</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t>, Total) y CAJA</w:t></w:r>' +
      '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
      '<w:r><w:t>(</w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>NoCaja</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t>)</w:t></w:r>' +
      '<w:r><w:t xml:space="preserve"> (EX-01), </w:t></w:r>' +
      '<w:r><w:t xml:space="preserve">muestra la ventana </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>Reporte</w:t></w:r>' +
      '<w:r><w:t>View</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
      '<w:r><w:t xml:space="preserve">con </w:t></w:r>' +
      '<w:r><w:t>una tabla con la información obtenida, un botón &#8220;Imprimir&#8221; y un botón &#8220;Regresar&#8221;.</w:t></w:r>' +
      '</w:p>'

# --- Paragraph: "El sistema envía al cuadro de impresión ..." (Flujo normal, step 3)
$p2 = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
      '<w:r><w:t xml:space="preserve">El sistema envía al cuadro de impresión predeterminada del sistema Windows y cierra la ventana </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>ReporteView</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t>.</w:t></w:r>' +
      '</w:p>'

# --- Paragraph: "El sistema cierra la ventana ReporteVentaView." (Flujo alterno, step 2)
$p3 = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
      '<w:r><w:t>El sistema cierra la ventana</w:t></w:r>' +
      '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>ReporteView</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t>.</w:t></w:r>' +
      '</w:p>'

# --- Paragraph: "El sistema muestra en pantalla la ventana ErrorView ..." (Excepciones, step 1)
$p4 = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr>' +
      '<w:r><w:t xml:space="preserve">El sistema muestra en pantalla la ventana </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>ErrorView</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t xml:space="preserve"> con el mensaje &#8220;No se pudo conectar a </w:t></w:r>' +
      '<w:r><w:t>la red de la empresa</w:t></w:r>' +
      '<w:r><w:t>, por favor revise su conexión&#8221; junto con un botón de aceptar.</w:t></w:r>' +
      '</w:p>'

# --- Paragraph: "El sistema cierra las ventanas ErrorView y ReporteVentaView." (Excepciones, step 4)
$p5 = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr>' +
      '<w:r><w:t xml:space="preserve">El sistema cierra las ventanas </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>ErrorView</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t xml:space="preserve"> y </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>ReporteView</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:t xml:space="preserve">. </w:t></w:r>' +
      '</w:p>'

Set-ParagraphXml 20 $p1
Set-ParagraphXml 22 $p2
Set-ParagraphXml 28 $p3
Set-ParagraphXml 33 $p4
Set-ParagraphXml 35 $p5
